$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C. This shifts the existing C:F data
# (Utilizador1..Utilizador4 results) one column to the right, becoming D:G.
$ws.Columns("C").Insert()

# Fill in the answers for the new respondent in column C (rows 2-11,
# one value per SUS question), plus the computed SUS score in row 12.
$ws.Range("C2").Value = 4
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 5
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 4
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 4
$ws.Range("C9").Value = 2
$ws.Range("C10").Value = 4
$ws.Range("C11").Value = 1
$ws.Range("C12").Value = 87.5

# Match the visual formatting of the other "SUS Score" row cells.
$ws.Range("C12").Interior.Color = $ws.Range("D12").Interior.Color
$ws.Range("C12").Borders.LineStyle = $ws.Range("D12").Borders.LineStyle
$ws.Range("C12").Font.Bold = $false

# The new column is a bit wider than the rest (not an auto best-fit column).
$ws.Columns("C").ColumnWidth = 11.9

# Move the active selection, as left by the editor, to C4.
$ws.Range("C4").Select()

Write-Host "done"
